$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new product row (SKU 3614222915713-0) under the existing table ---
$ws.Cells.Item(8,1).Value = "3614222915713-0"
$ws.Cells.Item(8,2).Value = 1
$ws.Cells.Item(8,3).Formula = "=B8/1.3"
$ws.Cells.Item(8,4).Value = 2

# --- Header row (SKU / Stock Quantity) -> smaller bold Arial ---
$ws.Range("A1").Font.Name = "Arial"
$ws.Range("A1").Font.Size = 10
$ws.Range("A1").Font.Bold = $true
$ws.Range("D1").Font.Name = "Arial"
$ws.Range("D1").Font.Size = 10
$ws.Range("D1").Font.Bold = $true

# --- Header row (Regular Price / Sale Price) -> bold Calibri 11, Currency style ---
$ws.Range("B1:C1").Style = "Currency"
$ws.Range("B1:C1").Font.Bold = $true
$ws.Range("B1:C1").Font.Size = 11

# --- Body rows: SKU / Stock Quantity columns revert to plain Arial 10 / General ---
$ws.Range("A2:A8").Font.Name = "Arial"
$ws.Range("A2:A8").Font.Size = 10
$ws.Range("A2:A8").Font.Bold = $false
$ws.Range("A2:A8").NumberFormat = "General"

$ws.Range("D2:D8").Font.Name = "Arial"
$ws.Range("D2:D8").Font.Size = 10
$ws.Range("D2:D8").Font.Bold = $false
$ws.Range("D2:D8").NumberFormat = "General"

# --- Body rows: Regular Price / Sale Price columns -> Currency cell style ---
$ws.Range("B2:C8").Style = "Currency"

# --- Row heights tightened to match the new (smaller) font metrics ---
for ($r = 1; $r -le 8; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.8
}

# --- Selection moved back to the header row ---
$ws.Range("A1:D1").Select()
$ws.Range("A1").Activate()
